# "Actualización modelo información planilla"
#
# The table "Tabla1" on the sheet had an AutoFilter criterion applied to its
# first column (NIT), showing only the row whose NIT equals 899999034 and
# hiding every other data row. This script clears that filter criterion
# (Excel's "Clear Filter" for the column), which unhides all the
# previously-filtered data rows while keeping the table/AutoFilter itself
# (and its drop-down arrows) in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = $ws.ListObjects.Item("Tabla1")

# Clear the AutoFilter criterion on column 1 (NIT) - same as picking
# "Clear Filter From NIT" in the Excel UI. This recomputes which rows are
# hidden, unhiding all rows that were only hidden due to that filter.
$table.Range.AutoFilter(1)
